$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 73.38544233333333
$ws.Range("H2").Value = 220.156327
$ws.Range("I2").Value = 0.1214979676060253
$ws.Range("J2").Value = 0.1214979676060253
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 462.7572979958807
$ws.Range("R2").Value = 4164.815681962926
$ws.Range("S2").Value = 0.00164769530965859
$ws.Range("T2").Value = 0.00164769530965859
$ws.Range("G3").Value = 73.38544233333333
$ws.Range("H3").Value = 220.156327
$ws.Range("I3").Value = 0.1214979676060253
$ws.Range("J3").Value = 0.1214979676060253
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 13395.16824787915
$ws.Range("R3").Value = 120556.5142309124
$ws.Range("S3").Value = 0.0476948845317068
$ws.Range("T3").Value = 0.04769488453170678
$ws.Range("G4").Value = 73.38544233333333
$ws.Range("H4").Value = 220.156327
$ws.Range("I4").Value = 0.1214979676060253
$ws.Range("J4").Value = 0.1214979676060253
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 9349.040651976171
$ws.Range("R4").Value = 84141.36586778554
$ws.Range("S4").Value = 0.03328822797345868
$ws.Range("T4").Value = 0.03328822797345867
$ws.Range("G5").Value = 73.38544233333333
$ws.Range("H5").Value = 220.156327
$ws.Range("I5").Value = 0.1214979676060253
$ws.Range("J5").Value = 0.1214979676060253
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 1425.438929729923
$ws.Range("R5").Value = 12828.95036756931
$ws.Range("S5").Value = 0.005075423010922808
$ws.Range("T5").Value = 0.005075423010922806
$ws.Range("G6").Value = 73.38544233333333
$ws.Range("H6").Value = 220.156327
$ws.Range("I6").Value = 0.1214979676060253
$ws.Range("J6").Value = 0.1214979676060253
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 9490.451733014766
$ws.Range("R6").Value = 85414.0655971329
$ws.Range("S6").Value = 0.03379173678027845
$ws.Range("T6").Value = 0.03379173678027843
$ws.Range("I7").Value = 0.3924995450689984
$ws.Range("J7").Value = 0.3924995450689983
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 1494.938825065043
$ws.Range("R7").Value = 13454.44942558539
$ws.Range("S7").Value = 0.005322884589727469
$ws.Range("T7").Value = 0.005322884589727466
$ws.Range("I8").Value = 0.3924995450689984
$ws.Range("J8").Value = 0.3924995450689983
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.1540784660819706
$ws.Range("T8").Value = 0.1540784660819705
$ws.Range("I9").Value = 0.3924995450689984
$ws.Range("J9").Value = 0.3924995450689983
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 30202.10358276185
$ws.Range("R9").Value = 271818.9322448566
$ws.Range("S9").Value = 0.1075377192983407
$ws.Range("T9").Value = 0.1075377192983406
$ws.Range("I10").Value = 0.3924995450689984
$ws.Range("J10").Value = 0.3924995450689983
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 4604.88469450652
$ws.Range("R10").Value = 41443.96225055867
$ws.Range("S10").Value = 0.0163961691053105
$ws.Range("T10").Value = 0.0163961691053105
$ws.Range("I11").Value = 0.3924995450689984
$ws.Range("J11").Value = 0.3924995450689983
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 30658.9324998952
$ws.Range("R11").Value = 275930.3924990568
$ws.Range("S11").Value = 0.1091643059936493
$ws.Range("T11").Value = 0.1091643059936492
$ws.Range("G12").Value = 138.1628113333333
$ws.Range("H12").Value = 414.488434
$ws.Range("I12").Value = 0.2287442882675098
$ws.Range("J12").Value = 0.2287442882675098
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 871.2334111950547
$ws.Range("R12").Value = 7841.100700755492
$ws.Range("S12").Value = 0.003102116836322103
$ws.Range("T12").Value = 0.003102116836322102
$ws.Range("G13").Value = 138.1628113333333
$ws.Range("H13").Value = 414.488434
$ws.Range("I13").Value = 0.2287442882675098
$ws.Range("J13").Value = 0.2287442882675098
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 25219.09038857627
$ws.Range("R13").Value = 226971.8134971864
$ws.Range("S13").Value = 0.08979518448887447
$ws.Range("T13").Value = 0.08979518448887443
$ws.Range("G14").Value = 138.1628113333333
$ws.Range("H14").Value = 414.488434
$ws.Range("I14").Value = 0.2287442882675098
$ws.Range("J14").Value = 0.2287442882675098
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 17601.44381060619
$ws.Range("R14").Value = 158412.9942954557
$ws.Range("S14").Value = 0.06267176452009886
$ws.Range("T14").Value = 0.06267176452009884
$ws.Range("G15").Value = 138.1628113333333
$ws.Range("H15").Value = 414.488434
$ws.Range("I15").Value = 0.2287442882675098
$ws.Range("J15").Value = 0.2287442882675098
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 2683.674631555748
$ws.Range("R15").Value = 24153.07168400174
$ws.Range("S15").Value = 0.009555501603571718
$ws.Range("T15").Value = 0.009555501603571714
$ws.Range("G16").Value = 138.1628113333333
$ws.Range("H16").Value = 414.488434
$ws.Range("I16").Value = 0.2287442882675098
$ws.Range("J16").Value = 0.2287442882675098
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 17867.67852813004
$ws.Range("R16").Value = 160809.1067531704
$ws.Range("S16").Value = 0.06361972081864273
$ws.Range("T16").Value = 0.0636197208186427
$ws.Range("G17").Value = 49.051656
$ws.Range("H17").Value = 147.154968
$ws.Range("I17").Value = 0.08121060965524597
$ws.Range("J17").Value = 0.08121060965524596
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 309.312188780976
$ws.Range("R17").Value = 2783.809699028784
$ws.Range("S17").Value = 0.001101338098571021
$ws.Range("T17").Value = 0.001101338098571021
$ws.Range("G18").Value = 49.051656
$ws.Range("H18").Value = 147.154968
$ws.Range("I18").Value = 0.08121060965524597
$ws.Range("J18").Value = 0.08121060965524596
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 8953.481290915945
$ws.Range("R18").Value = 80581.33161824349
$ws.Range("S18").Value = 0.03187979788119834
$ws.Range("T18").Value = 0.03187979788119833
$ws.Range("G19").Value = 49.051656
$ws.Range("H19").Value = 147.154968
$ws.Range("I19").Value = 0.08121060965524597
$ws.Range("J19").Value = 0.08121060965524596
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 6249.004045076808
$ws.Range("R19").Value = 56241.03640569127
$ws.Range("S19").Value = 0.0222502264139382
$ws.Range("T19").Value = 0.0222502264139382
$ws.Range("G20").Value = 49.051656
$ws.Range("H20").Value = 147.154968
$ws.Range("I20").Value = 0.08121060965524597
$ws.Range("J20").Value = 0.08121060965524596
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 952.779431546208
$ws.Range("R20").Value = 8575.014883915872
$ws.Range("S20").Value = 0.003392469891445861
$ws.Range("T20").Value = 0.00339246989144586
$ws.Range("G21").Value = 49.051656
$ws.Range("H21").Value = 147.154968
$ws.Range("I21").Value = 0.08121060965524597
$ws.Range("J21").Value = 0.08121060965524596
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 6343.524803978641
$ws.Range("R21").Value = 57091.72323580776
$ws.Range("S21").Value = 0.02258677737009256
$ws.Range("T21").Value = 0.02258677737009255
$ws.Range("G22").Value = 106.3337146666667
$ws.Range("H22").Value = 319.001144
$ws.Range("I22").Value = 0.1760475894022206
$ws.Range("J22").Value = 0.1760475894022206
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 670.5240292959413
$ws.Range("R22").Value = 6034.716263663472
$ws.Range("S22").Value = 0.002387470284896807
$ws.Range("T22").Value = 0.002387470284896806
$ws.Range("G23").Value = 106.3337146666667
$ws.Range("H23").Value = 319.001144
$ws.Range("I23").Value = 0.1760475894022206
$ws.Range("J23").Value = 0.1760475894022206
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 19409.27182685931
$ws.Range("R23").Value = 174683.4464417338
$ws.Range("S23").Value = 0.06910872349611089
$ws.Range("T23").Value = 0.06910872349611087
$ws.Range("G24").Value = 106.3337146666667
$ws.Range("H24").Value = 319.001144
$ws.Range("I24").Value = 0.1760475894022206
$ws.Range("J24").Value = 0.1760475894022206
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 13546.53170282453
$ws.Range("R24").Value = 121918.7853254208
$ws.Range("S24").Value = 0.04823382979707014
$ws.Range("T24").Value = 0.04823382979707014
$ws.Range("G25").Value = 106.3337146666667
$ws.Range("H25").Value = 319.001144
$ws.Range("I25").Value = 0.1760475894022206
$ws.Range("J25").Value = 0.1760475894022206
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 2065.426215463619
$ws.Range("R25").Value = 18588.83593917258
$ws.Range("S25").Value = 0.007354164056199485
$ws.Range("T25").Value = 0.007354164056199482
$ws.Range("G26").Value = 106.3337146666667
$ws.Range("H26").Value = 319.001144
$ws.Range("I26").Value = 0.1760475894022206
$ws.Range("J26").Value = 0.1760475894022206
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 13751.4329075289
$ws.Range("R26").Value = 123762.8961677601
$ws.Range("S26").Value = 0.04896340176794329
$ws.Range("T26").Value = 0.04896340176794328
